# Apply the edit: add the "Abschluss" checkmark (✔️) to column A for rows 2-5,
# matching the value already used in rows 6-10, and set the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use Value2 to avoid the Value property's automatic (and here erroneous)
# date/number coercion of the checkmark + variation-selector string.
$check = $ws.Range("A6").Value2

$ws.Range("A2").Value2 = $check
$ws.Range("A3").Value2 = $check
$ws.Range("A4").Value2 = $check
$ws.Range("A5").Value2 = $check

$ws.Range("A2").Select()
